$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 data values (new day of stats: nation.xlsx COVID dashboard refresh) ---
$ws.Range("A2").Value = 43916
$ws.Range("B2").Value = 11658
$ws.Range("C2").Value = 2129
$ws.Range("D2").Value = 578
$ws.Range("E2").Value = 9782
$ws.Range("F2").Value = 894
$ws.Range("G2").Value = 741
$ws.Range("H2").Value = 241

# --- "update death data source": TotalUKDeaths (C2, NewUKCases col actually carries the
# death-source format swap) now uses the "???,??0" number format (same family as the old
# EnglandCases style) instead of the filled "###,##0" style ---
$ws.Range("C2").NumberFormat = "???,??0"

# --- update nIreland / general refresh: move the active selection ---
[void]$ws.Range("B8").Select()
